$wb = $excel.ActiveWorkbook

# --- Cardiac sheet ---
$ws = $wb.Worksheets.Item("Cardiac")

# Row 2 / Row 3 content swap (retinker):
#   old row2: answer=FALSE, info="Pain not worse with exertion (requires they clarify exercise 1hr after meal)"
#   old row3: answer=blank, info="Do you have any PMHx? (counts as 2 independent minor features)"
#   new row2: answer=blank, info="Do you have any PMHx? (counts as 2 independent minor features)"
#   new row3: answer=TRUE,  info="Pain not worse with exertion (requires they clarify exercise 1hr after meal)"
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Do you have any PMHx? (counts as 2 independent minor features)"
$ws.Range("A3").Value = $true
$ws.Range("B3").Value = "Pain not worse with exertion (requires they clarify exercise 1hr after meal)"

# Row 9: remove stray space before the opening parenthesis
$ws.Range("B9").Value = "Alternative cause of esoph dysphagia becomes obvious(food gets stuck or relieved by regurgitation of food)"

# --- GERD sheet: add LR estimation (mark heartburn as a feature) ---
$ws = $wb.Worksheets.Item("GERD")
$ws.Range("A2").Value = $true

# --- Esophageal Dysphagia sheet: add LR estimation ---
$ws = $wb.Worksheets.Item("Esophageal Dysphagia")
$ws.Range("A2").Value = $true
$ws.Range("A7").Value = $true

# --- RA sheet: add LR estimation ---
$ws = $wb.Worksheets.Item("RA")
$ws.Range("A2").Value = $true

# --- CREST sheet: add LR estimation ---
$ws = $wb.Worksheets.Item("CREST")
$ws.Range("A5").ClearContents()
$ws.Range("A6").Value = $true
